# Append the 2025-05-11 bitcoin buy as a new row at the bottom of the
# existing data table (rows 1..13 -> new row 14).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds the date as literal text (matching the existing rows,
# e.g. "04/27/2025", "05/04/2025", "05/07/2025"), not a date serial.
# A leading apostrophe forces Excel to store it as text instead of
# auto-converting the "mm/dd/yyyy"-looking string into a date value.
$ws.Cells.Item(14, 1).Value = "'05/11/2025"
# Reset to the default "Normal" style so the cell carries no explicit
# number format override (consistent with the other text-date cells).
$ws.Cells.Item(14, 1).Style = "Normal"

$ws.Cells.Item(14, 2).Value = 0.0004785600000000003
$ws.Cells.Item(14, 3).Value = 104480.1069876295
$ws.Cells.Item(14, 4).Value = 50
